# Uppercase the "Drink" / "Eat" / "Smoke" category labels in column A
# (Catégorie) to "DRINK" / "EAT" / "SMOKE", leaving every other cell
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $v = $cell.Value2
    if ($v -eq "Drink") {
        $cell.Value2 = "DRINK"
    } elseif ($v -eq "Eat") {
        $cell.Value2 = "EAT"
    } elseif ($v -eq "Smoke") {
        $cell.Value2 = "SMOKE"
    }
}
